$p = $ppt.ActivePresentation

# Slide 8's speaker notes describe where the ClientSideInstance.xml file
# lives; update it to describe the package-solution.json feature definition
# instead.
$s = $p.Slides.Item(8)
$notes = $s.NotesPage.Shapes.Item(2)
$notes.TextFrame.TextRange.Text = "The files referenced above are elements of a feature. The definition of the Feature is part of the ``solution`` object in the **./config/package-solution.json** file."
